$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed "cryptos" price/volume snapshot produced by the GitHub
# Actions scraper run. Column D ("Price") holds numeric-looking values (e.g.
# "61.906.10", "0.0000230") that must stay TEXT, exactly like the rest of the
# sheet (which stores every data cell as a string). Writing them with a
# leading apostrophe forces a text/quoted entry -- the same thing Excel does
# when a user types an apostrophe before a number -- so the value keeps its
# exact textual form instead of being parsed into a number/scientific notation.

# row 2: 'Bitcoin'
$ws.Range("D2").Value2 = '''61.906.10'
$ws.Range("E2").Value2 = '  -1.04%  '

# row 3: 'Ethereum'
$ws.Range("D3").Value2 = '''2.895.92'
$ws.Range("E3").Value2 = '  -1.84%  '

# row 4: 'TetherUSD'
$ws.Range("E4").Value2 = '  +0.04%  '

# row 5: 'BNB'
$ws.Range("D5").Value2 = '''568.35'
$ws.Range("E5").Value2 = '  -3.51%  '

# row 6: 'Solana'
$ws.Range("D6").Value2 = '''143.86'
$ws.Range("E6").Value2 = '  -1.75%  '

# row 7: 'USDC'
$ws.Range("E7").Value2 = '  -0.09%  '

# row 8: 'XRP'
$ws.Range("E8").Value2 = '  -1.01%  '

# row 9: 'LidoStakedEther'
$ws.Range("D9").Value2 = '''2.894.81'
$ws.Range("E9").Value2 = '  -1.84%  '

# row 10: 'Toncoin'
$ws.Range("D10").Value2 = '''6.91'
$ws.Range("E10").Value2 = '  -0.78%  '

# row 11: 'Dogecoin'
$ws.Range("D11").Value2 = '''0.146'
$ws.Range("E11").Value2 = '  -2.43%  '

# row 12: 'Cardano'
$ws.Range("E12").Value2 = '  -1.24%  '

# row 13: 'ShibaInu'
$ws.Range("D13").Value2 = '''0.0000230'
$ws.Range("E13").Value2 = '  -1.33%  '

# row 14: 'Avalanche'
$ws.Range("D14").Value2 = '''32.17'
$ws.Range("E14").Value2 = '  -0.34%  '

# row 15: 'TRON'
$ws.Range("E15").Value2 = '  -0.07%  '

# row 16: 'WrappedliquidstakedEther2.0'
$ws.Range("D16").Value2 = '''3.375.36'
$ws.Range("E16").Value2 = '  -1.89%  '

# row 17: 'WrappedBTC'
$ws.Range("D17").Value2 = '''61.815.13'
$ws.Range("E17").Value2 = '  -1.17%  '

# row 18: 'WrappedEther' -> 'Polkadot'
$ws.Range("B18").Value2 = 'Polkadot'
$ws.Range("C18").Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value2 = '''6.52'
$ws.Range("E18").Value2 = '  -1.89%  '

# row 19: 'Polkadot' -> 'WrappedEther'
$ws.Range("B19").Value2 = 'WrappedEther'
$ws.Range("C19").Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value2 = '''2.889.35'
$ws.Range("E19").Value2 = '  -2.17%  '

# row 20: 'BitcoinCash'
$ws.Range("D20").Value2 = '''429.99'
$ws.Range("E20").Value2 = '  -1.19%  '

# row 21: 'Chainlink'
$ws.Range("D21").Value2 = '''12.95'
$ws.Range("E21").Value2 = '  -3.63%  '

# row 22: 'Polygon'
$ws.Range("D22").Value2 = '''0.653'
$ws.Range("E22").Value2 = '  -1.41%  '

# row 23: 'Uniswap'
$ws.Range("E23").Value2 = '  -1.44%  '

# row 24: 'Litecoin'
$ws.Range("D24").Value2 = '''78.86'
$ws.Range("E24").Value2 = '  -1.51%  '

# row 25: 'InternetComputer(DFINITY)'
$ws.Range("D25").Value2 = '''12.06'
$ws.Range("E25").Value2 = '  +1.33%  '

# row 26: 'RenderToken'
$ws.Range("E26").Value2 = '  -9.63%  '

# row 27: 'Dai'
$ws.Range("E27").Value2 = '  +0.03%  '

# row 28: 'Fetch.AI'
$ws.Range("D28").Value2 = '''2.02'
$ws.Range("E28").Value2 = '  -3.59%  '

# row 29: 'PEPE'
$ws.Range("D29").Value2 = '''0.0000110'
$ws.Range("E29").Value2 = '  +8.53%  '

# row 30: 'NEARProtocol'
$ws.Range("D30").Value2 = '''6.98'
$ws.Range("E30").Value2 = '  -3.22%  '

# row 31: 'PancakeSwap'
$ws.Range("E31").Value2 = '  -3.20%  '

# row 32: 'ImmutableX'
$ws.Range("E32").Value2 = '  -6.12%  '

# row 33: 'Hedera' -> 'FirstDigitalUSD'
$ws.Range("B33").Value2 = 'FirstDigitalUSD'
$ws.Range("C33").Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").Value2 = '''1.00'
$ws.Range("E33").Value2 = '  +0.01%  '

# row 34: 'FirstDigitalUSD' -> 'Hedera'
$ws.Range("B34").Value2 = 'Hedera'
$ws.Range("C34").Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value2 = '''0.107'
$ws.Range("E34").Value2 = '  -1.73%  '

# row 35: 'EthereumClassic'
$ws.Range("D35").Value2 = '''25.57'
$ws.Range("E35").Value2 = '  -2.42%  '

# row 36: 'Mantle'
$ws.Range("E36").Value2 = '  -3.87%  '

# row 37: 'Filecoin'
$ws.Range("E37").Value2 = '  -3.44%  '

# row 38: 'OKB'
$ws.Range("D38").Value2 = '''48.84'
$ws.Range("E38").Value2 = '  -1.55%  '

# row 39: 'dogwifhat'
$ws.Range("D39").Value2 = '''2.83'
$ws.Range("E39").Value2 = '  -6.50%  '

# row 40: 'Stacks'
$ws.Range("E40").Value2 = '  -5.06%  '

# row 41: 'Kaspa'
$ws.Range("E41").Value2 = '  +0.02%  '

# row 42: 'Cosmos'
$ws.Range("D42").Value2 = '''8.12'
$ws.Range("E42").Value2 = '  -2.74%  '

# row 43: 'Arweave'
$ws.Range("D43").Value2 = '''40.20'
$ws.Range("E43").Value2 = '  +2.62%  '

# row 44: 'TheGraph'
$ws.Range("D44").Value2 = '''0.268'
$ws.Range("E44").Value2 = '  -2.27%  '

# row 45: 'Maker'
$ws.Range("D45").Value2 = '''2.694.38'
$ws.Range("E45").Value2 = '  +0.50%  '

# row 46: 'VeChain'
$ws.Range("D46").Value2 = '''0.0335'
$ws.Range("E46").Value2 = '  -0.47%  '

# row 47: 'Monero'
$ws.Range("D47").Value2 = '''131.69'

# row 48: 'Bittensor'
$ws.Range("D48").Value2 = '''346.47'
$ws.Range("E48").Value2 = '  -2.71%  '

# row 50: 'Stellar'
$ws.Range("E50").Value2 = '  -1.40%  '

# row 51: 'InjectiveProtocol'
$ws.Range("D51").Value2 = '''21.62'
$ws.Range("E51").Value2 = '  -4.63%  '
